$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sprint-5")
$ws2 = $wb.Worksheets.Item("Test Report")

# ---- Sprint-5 sheet ----

# Header cell B2 ("Appetite") is cleared out.
$ws1.Range("B2").Value = $null

# Row 8/9/11: "Tested By" column now points at "Irvin" (shared-string index shift,
# same text as before).
$ws1.Range("G8").Value = "Irvin"
$ws1.Range("G9").Value = "Irvin"
$ws1.Range("G11").Value = "Irvin"

# Row 10 (TC3) gets new description/procedure/expected-output text, still tested by Irvin.
$ws1.Range("B10").Value = "To be able to view detailed information of the restaurant clicked"
$ws1.Range("C10").Value = "Implement backend yelp API call to retrieve pictures information"
$ws1.Range("D10").Value = "Appropriate pictures of each restaurant are displayed in detail."
$ws1.Range("G10").Value = "Irvin"

# Row 11 (TC4) gets swapped text (previously TC3's content).
$ws1.Range("B11").Value = "Call detail screen when you tap on the popup of map screen"
$ws1.Range("C11").Value = "Click on each popup and it should lead to the detailed screen of the restaurant.  "
$ws1.Range("D11").Value = "On clicking on map, popup appears, on further clicking on popup it shoulfd open the detail screen of the restaurant clicked."

# New rows 12 (TC5) and 13 (TC6), replacing the former blank filler rows 12-24.
$ws1.Range("A12").Value = "TC5"
$ws1.Range("B12").Value = "Restaurant's price category will be shown in the list of restaurants."
$ws1.Range("C12").Value = "Open the list of restaurants and check the information about price category."
$ws1.Range("D12").Value = "To be able to view price category for a choose restaurant in the list screen."
$ws1.Range("G12").Value = "Sergio Brunacci"
$ws1.Range("H12").NumberFormat = "@"
$ws1.Range("H12").Value = "03/25/2018"
$ws1.Range("H12").NumberFormat = "m/d/yy"
$ws1.Range("I12").Value = "Pass"
$ws1.Rows.Item(12).RowHeight = 43

$ws1.Range("A13").Value = "TC6"
$ws1.Range("B13").Value = "Restaurant's price category will be shown in the Details screen."
$ws1.Range("C13").Value = "Open the details screen and check the information about price category."
$ws1.Range("D13").Value = "To be able to view price category for a choose restaurant in the detail screen."
$ws1.Range("G13").Value = "Sergio Brunacci"
$ws1.Range("H13").NumberFormat = "@"
$ws1.Range("H13").Value = "03/25/2018"
$ws1.Range("H13").NumberFormat = "m/d/yy"
$ws1.Range("I13").Value = "Pass"
$ws1.Rows.Item(13).RowHeight = 43

# Drop the now-unused filler rows 14-24 entirely (shifts dimension/mergeCells down).
$ws1.Range("A14:J24").EntireRow.Delete()

# Update summary formulas to the new (smaller) data range.
$ws1.Range("B3").Formula = "=COUNTIF(I8:I12,""Pass"")"
$ws1.Range("D3").Formula = "=COUNTIF(I7:I669,""Pending"")"
$ws1.Range("D4").Formula = "=COUNTA(A8:A12)"

# ---- Test Report sheet ----
$ws2.Range("C8").Value = "Appetite"
$ws2.Activate()
$ws2.Range("C9").Select()

# ---- View state (Sprint-5 stays the active/selected tab) ----
$ws1.Activate()
$ws1.Range("B16").Select()
$excel.ActiveWindow.Zoom = 100
